# Generate Report for Handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values for the 940c288c-... row on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-13 06:44:50"
$wsZhCn.Range("H4").Value = "2016-03-13 06:45:09"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-13 06:44:54"
$wsDeDe.Range("H4").Value = "2016-03-13 06:45:16"
